# Regenerate the localization-status report for the handback: the
# handback files are now in sync with en-US, so flip the "Ready for
# handoff" status to "Handed back: in sync with en-US", refresh the
# handback timestamps, and clear the stale "version not latest" error.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-locale status columns (zh-cn / de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: Status, Latest Handback DateTime, Error Detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-16 18:43:31"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: Status, Latest Handback DateTime, Error Detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-16 18:43:38"
$wsDeDe.Range("P2").Value = ""

# --- Column widths: the longer status text (and the now-empty Error
# Detail column) changes the autofit report column widths. ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
